$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 09:22 AM"

# --- distance from Dma50 sheet: update stock order/values ---
$ws = $wb.Worksheets.Item("distance from Dma50")

$ws.Range("C2").Value = 9.5312
$ws.Range("C3").Value = 7.6947
$ws.Range("C4").Value = 5.7615
$ws.Range("C5").Value = 5.7077
$ws.Range("B6").Value = "NIFTYCOMMODITIES"
$ws.Range("C6").Value = 5.1975
$ws.Range("B7").Value = "CNXINFRA"
$ws.Range("C7").Value = 5.1928
$ws.Range("C8").Value = 4.654
$ws.Range("C9").Value = 4.5757
$ws.Range("C10").Value = 3.7626
$ws.Range("C11").Value = 3.6689
$ws.Range("C12").Value = 3.4202
$ws.Range("B13").Value = "NIFTY"
$ws.Range("C13").Value = 3.303
$ws.Range("B14").Value = "CNXENERGY"
$ws.Range("C14").Value = 3.3007
$ws.Range("C15").Value = 3.2214
$ws.Range("C16").Value = 3.1713
$ws.Range("C17").Value = 3.0063
$ws.Range("C18").Value = 2.9268
$ws.Range("B19").Value = "NIFTYCPSE"
$ws.Range("C19").Value = 2.5661
$ws.Range("B20").Value = "CNXNIFTYJUNIOR"
$ws.Range("C20").Value = 2.5014
$ws.Range("C21").Value = 2.4051
$ws.Range("C22").Value = 1.7927
$ws.Range("C23").Value = 1.5128
$ws.Range("C24").Value = 1.4446
$ws.Range("C25").Value = 1.2717
$ws.Range("C26").Value = 1.1876
$ws.Range("C27").Value = 0.7259
$ws.Range("C28").Value = 0.6994
$ws.Range("C29").Value = 0.2416
$ws.Range("C30").Value = -1.9081
